$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: columns G6:AJ6 were empty inline-string cells; fill with numeric odds values
$ws.Range("G6").Value = 2.15
$ws.Range("H6").Value = 2.85
$ws.Range("I6").Value = 3.6
$ws.Range("J6").Value = 1.1
$ws.Range("K6").Value = 5.8
$ws.Range("L6").Value = 1.47
$ws.Range("M6").Value = 2.35
$ws.Range("N6").Value = 2.35
$ws.Range("O6").Value = 1.47
$ws.Range("P6").Value = 1.5
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.65
$ws.Range("T6").Value = 5.7
$ws.Range("U6").Value = 9.25
$ws.Range("V6").Value = 9.25
$ws.Range("W6").Value = 21
$ws.Range("X6").Value = 21
$ws.Range("Y6").Value = 40
$ws.Range("Z6").Value = 6.4
$ws.Range("AA6").Value = 5.7
$ws.Range("AB6").Value = 17
$ws.Range("AC6").Value = 100
$ws.Range("AD6").Value = 101
$ws.Range("AE6").Value = 8.25
$ws.Range("AF6").Value = 18
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 55
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 55

# Row 10: K10 12 -> 13
$ws.Range("K10").Value = 13

# Row 11: several odds updated
$ws.Range("G11").Value = 2.1
$ws.Range("I11").Value = 3.1
$ws.Range("P11").Value = 1.29
$ws.Range("Q11").Value = 3.5
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 2.38
$ws.Range("W11").Value = 21
$ws.Range("AC11").Value = 34
$ws.Range("AG11").Value = 11

# Row 13: J13 1.02 -> 1.01, K13 21 -> 23
$ws.Range("J13").Value = 1.01
$ws.Range("K13").Value = 23
